# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# for every cell in column G that matches exactly, across the whole used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G is the 7th column
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
